# Apply the "Inventory by Room" -> "Inventory" rename and related view/selection
# tidy-up described by the commit "Changed invetory by room to just inventory".

$wb = $excel.ActiveWorkbook

# 1) Rename the "Inventory by Room" sheet to "Inventory".
#    (This also re-points the _xlnm._FilterDatabase defined name that referenced
#    'Inventory by Room'!$A$1:$C$1 so it now reads Inventory!$A$1:$C$1.)
$wsInventory = $wb.Worksheets.Item("Inventory by Room")
$wsInventory.Name = "Inventory"

# 2) Remove the 4 blank formatting-only rows (9:12) from "Event Requirements",
#    shifting everything below up and shrinking the used range from H129 to H125.
$wsEventReq = $wb.Worksheets.Item("Event Requirements")
$wsEventReq.Rows("9:12").Delete()

# 3) Update the lingering selection on "Event Requirements (2)" to A8:H12.
$wsEventReq2 = $wb.Worksheets.Item("Event Requirements (2)")
[void]$wsEventReq2.Range("A8:H12").Select()

# 4) Update the selection on "Event Requirements" to the new A8:XFD11 block.
$wsEventReq.Activate()
[void]$wsEventReq.Range("A8:XFD11").Select()

# 5) Finally, make "Inventory" the active/visible tab when the workbook is reopened.
$wsInventory.Activate()
